$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1113.6666
$ws.Range("I127").Value = 836.8
$ws.Range("K127").Value = 2510.4
$ws.Range("M127").Value = 2449.6
$ws.Range("H131").Value = 8628.708000000001
$ws.Range("I131").Value = 1094.6666
$ws.Range("K131").Value = 3283.9998
$ws.Range("M131").Value = 1756.0002
$ws.Range("H135").Value = 3206.6365
$ws.Range("I135").Value = 3168.8572
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 28519.7148
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -25984.7148
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 591940.5
$ws.Range("I137").Value = 1003113.1
$ws.Range("J137").Value = 4551.143
$ws.Range("K137").Value = 3009339.3
$ws.Range("L137").Value = 13653.429
$ws.Range("M137").Value = -3006789.3
$ws.Range("N137").Value = -18753.429
$ws.Range("H138").Value = 554888.1
$ws.Range("J138").Value = 7053.8335
$ws.Range("L138").Value = 21161.5005
$ws.Range("N138").Value = -31441.5005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 29413.428
$ws.Range("I61").Value = 33482.332
$ws.Range("K61").Value = 33482.332
$ws.Range("M61").Value = -33270.332
$ws.Range("H74").Value = 4215.1763
$ws.Range("I74").Value = 22578.834
$ws.Range("K74").Value = 22578.834
$ws.Range("M74").Value = -21704.834
$ws.Range("H77").Value = 4215.1763
$ws.Range("I77").Value = 22578.834
$ws.Range("K77").Value = 112894.17
$ws.Range("M77").Value = -108526.17
$ws.Range("H136").Value = 29413.428
$ws.Range("I136").Value = 33482.332
$ws.Range("K136").Value = 100446.996
$ws.Range("M136").Value = -97896.99600000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 64794.332
$ws.Range("I105").Value = 222779.8
$ws.Range("J105").Value = 4030.6924
$ws.Range("K105").Value = 222779.8
$ws.Range("L105").Value = 4030.6924
$ws.Range("M105").Value = -221032.8
$ws.Range("N105").Value = -7524.6924
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3524.6667
$ws.Range("I31").Value = 2633.7778
$ws.Range("J31").Value = 4861
$ws.Range("K31").Value = 2633.7778
$ws.Range("L31").Value = 4861
$ws.Range("M31").Value = -2338.7778
$ws.Range("N31").Value = -5451
$ws.Range("H34").Value = 3524.6667
$ws.Range("I34").Value = 2633.7778
$ws.Range("J34").Value = 4861
$ws.Range("K34").Value = 2633.7778
$ws.Range("L34").Value = 4861
$ws.Range("M34").Value = -2431.7778
$ws.Range("N34").Value = -5265
$ws.Range("H58").Value = 4549.077
$ws.Range("I58").Value = 4169
$ws.Range("J58").Value = 4786.625
$ws.Range("K58").Value = 4169
$ws.Range("L58").Value = 4786.625
$ws.Range("M58").Value = -3966
$ws.Range("N58").Value = -5192.625
$ws.Range("H99").Value = 12431716
$ws.Range("I99").Value = 18165738
$ws.Range("K99").Value = 18165738
$ws.Range("M99").Value = -18164240
$ws.Range("H126").Value = 12431716
$ws.Range("I126").Value = 18165738
$ws.Range("K126").Value = 54497214
$ws.Range("M126").Value = -54494744
$ws.Range("H132").Value = 9223.5
$ws.Range("I132").Value = 10729.909
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 32189.727
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -29659.727
$ws.Range("N132").Value = -16160
$ws.Range("H134").Value = 2220.2727
$ws.Range("I134").Value = 2114.7646
$ws.Range("J134").Value = 2579
$ws.Range("K134").Value = 6344.293799999999
$ws.Range("L134").Value = 7737
$ws.Range("M134").Value = -3809.293799999999
$ws.Range("N134").Value = -12807
$ws.Range("H136").Value = 4549.077
$ws.Range("I136").Value = 4169
$ws.Range("J136").Value = 4786.625
$ws.Range("K136").Value = 12507
$ws.Range("L136").Value = 14359.875
$ws.Range("M136").Value = -9957
$ws.Range("N136").Value = -19459.875
$ws.Range("H141").Value = 590139.0600000001
$ws.Range("J141").Value = 671828.9399999999
$ws.Range("L141").Value = 671828.9399999999
$ws.Range("N141").Value = -682188.9399999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 271802.4
$ws.Range("I5").Value = 857.0714
$ws.Range("J5").Value = 436725.66
$ws.Range("K5").Value = 2571.2142
$ws.Range("L5").Value = 1310176.98
$ws.Range("M5").Value = -2459.2142
$ws.Range("N5").Value = -1310400.98
$ws.Range("H59").Value = 5200
$ws.Range("J59").Value = 6000
$ws.Range("L59").Value = 18000
$ws.Range("N59").Value = -19080
$ws.Range("H113").Value = 2182
$ws.Range("I113").Value = 2000.4286
$ws.Range("J113").Value = 2363.5715
$ws.Range("K113").Value = 6001.2858
$ws.Range("L113").Value = 7090.7145
$ws.Range("M113").Value = -3831.2858
$ws.Range("N113").Value = -11430.7145
$ws.Range("H122").Value = 5766.1177
$ws.Range("J122").Value = 7114.4614
$ws.Range("L122").Value = 64030.1526
$ws.Range("N122").Value = -68930.1526
$ws.Range("H135").Value = 271802.4
$ws.Range("I135").Value = 857.0714
$ws.Range("J135").Value = 436725.66
$ws.Range("K135").Value = 7713.6426
$ws.Range("L135").Value = 3930530.94
$ws.Range("M135").Value = -5178.6426
$ws.Range("N135").Value = -3935600.94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4792.952
$ws.Range("I16").Value = 4787.75
$ws.Range("J16").Value = 4809.6
$ws.Range("K16").Value = 4787.75
$ws.Range("L16").Value = 4809.6
$ws.Range("M16").Value = -4617.75
$ws.Range("N16").Value = -5149.6
$ws.Range("H68").Value = 3658.5715
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 4152.5
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 4152.5
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -5650.5
$ws.Range("H71").Value = 3658.5715
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 4152.5
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 20762.5
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -28250.5
$ws.Range("H136").Value = 10680.1875
$ws.Range("I136").Value = 11398.3
$ws.Range("J136").Value = 9483.333000000001
$ws.Range("K136").Value = 34194.89999999999
$ws.Range("L136").Value = 28449.999
$ws.Range("M136").Value = -31644.89999999999
$ws.Range("N136").Value = -33549.999
$ws.Range("H140").Value = 179000
$ws.Range("J140").Value = 179000
$ws.Range("L140").Value = 179000
$ws.Range("N140").Value = -189360
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H132").Value = 11062.54
$ws.Range("I132").Value = 11655
$ws.Range("K132").Value = 34965
$ws.Range("M132").Value = -32435
$ws.Range("H136").Value = 2208434.2
$ws.Range("I136").Value = 5143021.5
$ws.Range("J136").Value = 7493.75
$ws.Range("K136").Value = 15429064.5
$ws.Range("L136").Value = 22481.25
$ws.Range("M136").Value = -15426514.5
$ws.Range("N136").Value = -27581.25
